$d = $word.ActiveDocument

# --- "First game" (pictures) positive-framing question: drop "my favorite" ---
$d.Content.Find.Execute(
    "Can you tell me which one of these pictures is LIKE my favorite picture on the top?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Can you tell me which one of these pictures is LIKE the picture on the top?",
    2) | Out-Null

# --- "First game" (pictures) negative-framing question: drop "my favorite" ---
$d.Content.Find.Execute(
    "Can you tell me which one of these pictures is NOT LIKE my favorite picture on the top?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Can you tell me which one of these pictures is NOT LIKE the picture on the top?",
    2) | Out-Null

# --- Familiarization trial narration: drop the leading "That was fun! " ---
$d.Content.Find.Execute(
    "That was fun! Let" + [char]8217 + "s watch it one more time!",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Let" + [char]8217 + "s watch it one more time!",
    2) | Out-Null
